$wb = $excel.ActiveWorkbook

# --- Bug List sheet: add four new bug/enhancement rows (8-11) ---
$bugList = $wb.Worksheets.Item("Bug List")

$bugList.Range("A8").Value = "Add in the trailer name to the Paired with trailer screen not just the rego number"
$bugList.Range("B8").Value = "1.0.6"

$bugList.Range("A9").Value = "Display the Trailer pair and tdefault truck in the index page for trailers"
$bugList.Range("B9").Value = "1.0.6"

$bugList.Range("A10").Value = "Customer orer sheet displays XX in truck type"
$bugList.Range("B10").Value = "1.0.6"

$bugList.Range("A11").Value = "Customer order sheet to be emailed on submit"
$bugList.Range("B11").Value = "1.0.6"

# --- Make "Bug List" the active sheet/tab with B11 selected ---
$bugList.Activate()
$bugList.Range("B11").Select()
